# Applies the "Update dashboards - 2025-12-20" data refresh to the
# Economic Dashboard worksheet (rows 28-30, 42-43, 47-52).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28 (Mich NTM Inflation Exp / UMCSENT) ---
$ws.Range("N28").Value = 45962
$ws.Range("Q28").Value = 51
$ws.Range("R28").Value = 53.6
$ws.Range("S28").Value = 55.1
$ws.Range("T28").Value = 58.2
$ws.Range("U28").Value = 61.7

# --- Row 29 (5yr, 5yr Forward / T5YIFR) ---
$ws.Range("N29").Value = 46010
$ws.Range("Q29").Value = 2.21
$ws.Range("R29").Value = 2.22
$ws.Range("S29").Value = 2.22
$ws.Range("T29").Value = 2.21
$ws.Range("U29").Value = 2.21

# --- Row 30 (10yr TIPS / T10YIE) ---
$ws.Range("N30").Value = 46010
$ws.Range("Q30").Value = 2.24
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = 2.24
$ws.Range("T30").Value = 2.23
$ws.Range("U30").Value = 2.25

# --- Row 42 (Existing Home Sales / EXHOSLUSM495S) ---
$ws.Range("C42").Value = 45962
$ws.Range("F42").Value = 4130000
$ws.Range("G42").Value = 4110000
$ws.Range("H42").Value = 4050000
$ws.Range("I42").Value = 4000000
$ws.Range("J42").Value = 4010000

# --- Row 43 (Existing Home Sales Y/Y % Delta) ---
$ws.Range("C43").Value = 45962
$ws.Range("F43").Value = -0.009592326139088728

# --- Row 47 (FFR / DFF) ---
$ws.Range("N47").Value = 46009

# --- Row 48 (2y UST / DGS2) ---
$ws.Range("N48").Value = 46009
$ws.Range("Q48").Value = 3.46
$ws.Range("R48").Value = 3.49
$ws.Range("S48").Value = 3.48
$ws.Range("T48").Value = 3.51

# --- Row 49 (5y UST / DGS5) ---
$ws.Range("N49").Value = 46009
$ws.Range("Q49").Value = 3.66
$ws.Range("R49").Value = 3.7
$ws.Range("S49").Value = 3.69
$ws.Range("T49").Value = 3.73

# --- Row 50 (10y UST / DGS10) ---
$ws.Range("N50").Value = 46009
$ws.Range("Q50").Value = 4.12
$ws.Range("R50").Value = 4.16
$ws.Range("S50").Value = 4.15
$ws.Range("T50").Value = 4.18

# --- Row 52 (BAA / DBAA) ---
$ws.Range("N52").Value = 46009
$ws.Range("Q52").Value = 5.9
$ws.Range("R52").Value = 5.94
$ws.Range("S52").Value = 5.93
$ws.Range("T52").Value = 5.95
